$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("D1").EntireColumn.Insert()
Write-Host "Inserted"
Write-Host $ws.Range("E7").Value
